# ---------------------------------------------------------------------------
# Add a new worksheet "Filtering and Spin" with dilution / dry-cell-weight
# data, formatting, and a scatter chart, matching the target commit.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new worksheet after the last existing sheet -------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Filtering and Spin"

# --- 2. Header row (order chosen to match shared-string insertion order) ---
$ws.Range("A1").Value = "Date"
$ws.Range("E1").Value = "End Wt (g)"
$ws.Range("D1").Value = "Start Wt (g)"
$ws.Range("F1").Value = "Control Adjust"
$ws.Range("G1").Value = "Final Wt (g)"
$ws.Range("B6").Value = "Control 1"
$ws.Range("B7").Value = "Control 2"
$ws.Range("C6").Value = "N/A"
$ws.Range("C7").Value = "N/A"
$ws.Range("F6").Value = "Average Ctrl."
$ws.Range("G6").Value = "Ctrl. Stdev."
$ws.Range("I1").Value = "Final Density (g/L)"
$ws.Range("H1").Value = "Local Density (g/L)"
$ws.Range("B1").Value = "Sample #"
$ws.Range("C1").Value = "OD"

# --- 3. Data rows 2-5 (first batch of samples) ------------------------------
$ws.Range("A2").Value = 42341
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.63
$ws.Range("D2").Value = 3.096
$ws.Range("E2").Value = 3.119
$ws.Range("F2").Formula = "=E2-`$F`$7"
$ws.Range("G2").Formula = "=F2-D2"
$ws.Range("H2").Formula = "=G2/0.05"
$ws.Range("I2").Formula = "=H2/C2"

$ws.Range("A3").Value = 42341
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0.63
$ws.Range("D3").Value = 2.511
$ws.Range("E3").Value = 2.533
$ws.Range("F3").Formula = "=E3-`$F`$7"
$ws.Range("G3").Formula = "=F3-D3"
$ws.Range("H3").Formula = "=G3/0.05"
$ws.Range("I3").Formula = "=H3/C3"

$ws.Range("A4").Value = 42341
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0.626
$ws.Range("D4").Value = 2.566
$ws.Range("E4").Value = 2.593
$ws.Range("F4").Formula = "=E4-`$F`$7"
$ws.Range("G4").Formula = "=F4-D4"
$ws.Range("H4").Formula = "=G4/0.05"
$ws.Range("I4").Formula = "=H4/C4"

$ws.Range("A5").Value = 42341
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0.626
$ws.Range("D5").Value = 2.673
$ws.Range("E5").Value = 2.693
$ws.Range("F5").Formula = "=E5-`$F`$7"
$ws.Range("G5").Formula = "=F5-D5"
$ws.Range("H5").Formula = "=G5/0.05"
$ws.Range("I5").Formula = "=H5/C5"

# --- 4. Control rows 6-7 -----------------------------------------------------
$ws.Range("A6").Value = 42342
$ws.Range("D6").Value = 2.88
$ws.Range("E6").Value = 2.887

$ws.Range("A7").Value = 42342
$ws.Range("D7").Value = 2.602
$ws.Range("E7").Value = 2.61
$ws.Range("F7").Formula = "=AVERAGE((E6-D6),(E7-D7))"
$ws.Range("G7").Formula = "=STDEV((E6-D6),(E7-D7))"

# --- 5. Data rows 8-11 (second batch of samples) -----------------------------
$ws.Range("A8").Value = 42343
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.634
$ws.Range("D8").Value = 2.833
$ws.Range("E8").Value = 2.856
$ws.Range("F8").Formula = "=E8-`$F`$7"
$ws.Range("G8").Formula = "=F8-D8"
$ws.Range("H8").Formula = "=G8/0.05"
$ws.Range("I8").Formula = "=H8/C8"

$ws.Range("A9").Value = 42343
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 0.634
$ws.Range("D9").Value = 2.862
$ws.Range("E9").Value = 2.885
$ws.Range("F9").Formula = "=E9-`$F`$7"
$ws.Range("G9").Formula = "=F9-D9"
$ws.Range("H9").Formula = "=G9/0.05"
$ws.Range("I9").Formula = "=H9/C9"

$ws.Range("A10").Value = 42345
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0.662
$ws.Range("D10").Value = 3.002
$ws.Range("E10").Value = 3.028
$ws.Range("F10").Formula = "=E10-`$F`$7"
$ws.Range("G10").Formula = "=F10-D10"
$ws.Range("H10").Formula = "=G10/0.05"
$ws.Range("I10").Formula = "=H10/C10"

$ws.Range("A11").Value = 42345
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 0.662
$ws.Range("D11").Value = 3.148
$ws.Range("E11").Value = 3.173
$ws.Range("F11").Formula = "=E11-`$F`$7"
$ws.Range("G11").Formula = "=F11-D11"
$ws.Range("H11").Formula = "=G11/0.05"
$ws.Range("I11").Formula = "=H11/C11"

# --- 6. Trailing placeholder rows 12-15 -------------------------------------
$ws.Range("A12").Value = 42346
$ws.Range("B12").Value = 1

$ws.Range("A13").Value = 42346
$ws.Range("B13").Value = 2

$ws.Range("A14").Value = 42347
$ws.Range("B14").Value = 1

$ws.Range("A15").Value = 42347
$ws.Range("B15").Value = 2

# --- 7. Number formats (dates) ----------------------------------------------
$ws.Range("A2:A15").NumberFormat = "d-mmm"

# --- 8. Strikethrough formatting for superseded rows (2,3,5) ---------------
$ws.Range("A2:I2").Font.Strikethrough = $true
$ws.Range("A3:I3").Font.Strikethrough = $true
$ws.Range("A5:I5").Font.Strikethrough = $true

# --- 9. Black-out formatting for unused control summary cells --------------
$ws.Range("H6:I7").Interior.ThemeColor = 1

# --- 10. Column widths (best-fit columns as in the source) -----------------
$ws.Columns.Item(4).ColumnWidth = 10.21875
$ws.Columns.Item(6).ColumnWidth = 12.5546875
$ws.Columns.Item(7).ColumnWidth = 10.109375
$ws.Columns.Item(8).ColumnWidth = 16
$ws.Columns.Item(9).ColumnWidth = 15.44140625

# --- 11. Select a neutral cell like the source sheet (H17) -------------------
$ws.Range("H17").Select()

Write-Host "sheet built"
